$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header timestamp
$ws.Range("A1").Value = 'Datos actualizados a 14 de Octubre de 2020 a las 21:28'

# Update country data rows (values refreshed and table re-sorted by total cases)
$rows = @(
    @(4, 'Estados Unidos', 8126349, 36096, 5252194, 2652724, 0, 558, 221431),
    @(13, 'Francia', 779063, 22591, 103413, 642613, 0, 104, 33037),
    @(23, 'Alemania', 340633, 4954, 281900, 48972, 0, 21, 9761),
    @(24, 'Arabia Saudita', 340590, 501, 326820, 8662, 0, 21, 5108),
    @(25, 'Turquia', 340450, 1671, 298368, 33068, 0, 57, 9014),
    @(30, 'Canada', 188959, 2078, 159036, 20260, 0, 9, 9663),
    @(49, 'Costa Rica', 91780, 1542, 54155, 36491, 0, 10, 1134),
    @(50, 'Portugal', 91193, 2072, 54493, 34583, 0, 7, 2117),
    @(52, 'Etiopia', 86430, 712, 40165, 44953, 0, 7, 1312),
    @(90, 'Republica de Macedonia', 21636, 443, 16591, 4236, 0, 9, 809),
    @(91, 'Camerun', 21203, 0, 20117, 663, 0, 0, 423),
    @(107, 'Mozambique', 10392, 134, 8035, 2284, 0, 0, 73),
    @(108, 'Tayikistan', 10336, 39, 9245, 1012, 0, 0, 79),
    @(123, 'Malaui', 5829, 2, 4723, 925, 0, 0, 181),
    @(130, 'Congo', 5156, 38, 3887, 1177, 0, 2, 92),
    @(131, 'Trinidad yTobago', 5127, 0, 3367, 1667, 0, 0, 93),
    @(148, 'Mali', 3352, 55, 2550, 670, 0, 0, 132),
    @(166, 'Republica del Chad', 1329, 8, 1120, 117, 0, 0, 92)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
